# Update NATMI Rarres2-Ccrl2 LR-pairs sheet with newly recomputed TPM-based
# statistics (commit: "update scripts wuth new tpm").
# All values in this worksheet are static (no formulas), so the refreshed
# ligand/receptor/edge expression statistics are written directly per cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9919543333333333
$ws.Range("H2").Value = 2.975863
$ws.Range("I2").Value = 0.008811579445878926
$ws.Range("J2").Value = 0.008811579445878926
$ws.Range("M2").Value = 1.764388666666667
$ws.Range("N2").Value = 5.293166
$ws.Range("O2").Value = 0.2918165420774624
$ws.Range("P2").Value = 0.2918165420774624
$ws.Range("Q2").Value = 1.750192983584222
$ws.Range("R2").Value = 15.751736852258
$ws.Range("S2").Value = 0.002571364644137231
$ws.Range("T2").Value = 0.002571364644137231
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9919543333333333
$ws.Range("H3").Value = 2.975863
$ws.Range("I3").Value = 0.008811579445878926
$ws.Range("J3").Value = 0.008811579445878926
$ws.Range("O3").Value = 0.08474962093431622
$ws.Range("P3").Value = 0.08474962093431622
$ws.Range("Q3").Value = 0.5082926103664445
$ws.Range("R3").Value = 4.574633493298
$ws.Range("S3").Value = 0.0007467780178708511
$ws.Range("T3").Value = 0.0007467780178708511
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9919543333333333
$ws.Range("H4").Value = 2.975863
$ws.Range("I4").Value = 0.008811579445878926
$ws.Range("J4").Value = 0.008811579445878926
$ws.Range("M4").Value = 3.769421666666666
$ws.Range("N4").Value = 11.308265
$ws.Range("O4").Value = 0.6234338369882213
$ws.Range("P4").Value = 0.6234338369882212
$ws.Range("Q4").Value = 3.739094156410555
$ws.Range("R4").Value = 33.65184740769499
$ws.Range("S4").Value = 0.005493436783870843
$ws.Range("T4").Value = 0.005493436783870842
$ws.Range("I5").Value = 0.6711393126876655
$ws.Range("J5").Value = 0.6711393126876655
$ws.Range("M5").Value = 1.764388666666667
$ws.Range("N5").Value = 5.293166
$ws.Range("O5").Value = 0.2918165420774624
$ws.Range("P5").Value = 0.2918165420774624
$ws.Range("Q5").Value = 133.3045140531358
$ws.Range("R5").Value = 1199.740626478222
$ws.Range("S5").Value = 0.1958495534807594
$ws.Range("T5").Value = 0.1958495534807594
$ws.Range("I6").Value = 0.6711393126876655
$ws.Range("J6").Value = 0.6711393126876655
$ws.Range("O6").Value = 0.08474962093431622
$ws.Range("P6").Value = 0.08474962093431622
$ws.Range("S6").Value = 0.05687880234439717
$ws.Range("T6").Value = 0.05687880234439717
$ws.Range("I7").Value = 0.6711393126876655
$ws.Range("J7").Value = 0.6711393126876655
$ws.Range("M7").Value = 3.769421666666666
$ws.Range("N7").Value = 11.308265
$ws.Range("O7").Value = 0.6234338369882213
$ws.Range("P7").Value = 0.6234338369882212
$ws.Range("Q7").Value = 284.7903826573894
$ws.Range("R7").Value = 2563.113443916504
$ws.Range("S7").Value = 0.4184109568625089
$ws.Range("T7").Value = 0.4184109568625088
$ws.Range("G8").Value = 36.02919333333333
$ws.Range("H8").Value = 108.08758
$ws.Range("I8").Value = 0.3200491078664556
$ws.Range("J8").Value = 0.3200491078664556
$ws.Range("M8").Value = 1.764388666666667
$ws.Range("N8").Value = 5.293166
$ws.Range("O8").Value = 0.2918165420774624
$ws.Range("P8").Value = 0.2918165420774624
$ws.Range("Q8").Value = 63.56950038647556
$ws.Range("R8").Value = 572.1255034782801
$ws.Range("S8").Value = 0.09339562395256584
$ws.Range("T8").Value = 0.09339562395256584
$ws.Range("G9").Value = 36.02919333333333
$ws.Range("H9").Value = 108.08758
$ws.Range("I9").Value = 0.3200491078664556
$ws.Range("J9").Value = 0.3200491078664556
$ws.Range("O9").Value = 0.08474962093431622
$ws.Range("P9").Value = 0.08474962093431622
$ws.Range("Q9").Value = 18.46191111163111
$ws.Range("R9").Value = 166.15720000468
$ws.Range("S9").Value = 0.02712404057204819
$ws.Range("T9").Value = 0.02712404057204819
$ws.Range("G10").Value = 36.02919333333333
$ws.Range("H10").Value = 108.08758
$ws.Range("I10").Value = 0.3200491078664556
$ws.Range("J10").Value = 0.3200491078664556
$ws.Range("M10").Value = 3.769421666666666
$ws.Range("N10").Value = 11.308265
$ws.Range("O10").Value = 0.6234338369882213
$ws.Range("P10").Value = 0.6234338369882212
$ws.Range("Q10").Value = 135.8092219831889
$ws.Range("R10").Value = 1222.2829978487
$ws.Range("S10").Value = 0.1995294433418415
$ws.Range("T10").Value = 0.1995294433418415
